$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4800
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 5371.4287
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 5371.4287
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -5867.4287
$ws.Range("H67").Value = 4800
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 5371.4287
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 5371.4287
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -7087.4287
$ws.Range("H70").Value = 3611.6667
$ws.Range("I70").Value = 3570
$ws.Range("K70").Value = 10710
$ws.Range("M70").Value = -10440
$ws.Range("H73").Value = 3611.6667
$ws.Range("I73").Value = 3570
$ws.Range("K73").Value = 10710
$ws.Range("M73").Value = -9774
$ws.Range("H76").Value = 3383.3333
$ws.Range("I76").Value = 3422.2222
$ws.Range("K76").Value = 3422.2222
$ws.Range("M76").Value = -3107.2222
$ws.Range("H79").Value = 3383.3333
$ws.Range("I79").Value = 3422.2222
$ws.Range("K79").Value = 3422.2222
$ws.Range("M79").Value = -2330.2222
$ws.Range("H113").Value = 3852.7144
$ws.Range("I113").Value = 4774.5
$ws.Range("J113").Value = 3484
$ws.Range("K113").Value = 4774.5
$ws.Range("L113").Value = 3484
$ws.Range("M113").Value = -1520.5
$ws.Range("N113").Value = -9992
$ws.Range("H116").Value = 1802.4445
$ws.Range("I116").Value = 1584.3
$ws.Range("J116").Value = 2075.125
$ws.Range("K116").Value = 1584.3
$ws.Range("L116").Value = 2075.125
$ws.Range("M116").Value = 1857.7
$ws.Range("N116").Value = -8959.125
$ws.Range("H127").Value = 632.8333
$ws.Range("J127").Value = 1900
$ws.Range("L127").Value = 5700
$ws.Range("N127").Value = -15620
$ws.Range("H132").Value = 791954.9
$ws.Range("I132").Value = 1542.9056
$ws.Range("J132").Value = 5446603
$ws.Range("K132").Value = 4628.7168
$ws.Range("L132").Value = 16339809
$ws.Range("M132").Value = -2098.7168
$ws.Range("N132").Value = -16344869
$ws.Range("H138").Value = 3608996.2
$ws.Range("I138").Value = 479104.1
$ws.Range("J138").Value = 5956415.5
$ws.Range("K138").Value = 1437312.3
$ws.Range("L138").Value = 17869246.5
$ws.Range("M138").Value = -1432172.3
$ws.Range("N138").Value = -17879526.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 970.09
$ws.Range("I32").Value = 846.0112
$ws.Range("J32").Value = 1974
$ws.Range("K32").Value = 846.0112
$ws.Range("L32").Value = 1974
$ws.Range("M32").Value = -559.0112
$ws.Range("N32").Value = -2548
$ws.Range("H61").Value = 52737764
$ws.Range("I61").Value = 62563348
$ws.Range("K61").Value = 62563348
$ws.Range("M61").Value = -62563136
$ws.Range("H63").Value = 2850
$ws.Range("I63").Value = 2840
$ws.Range("J63").Value = 2900
$ws.Range("K63").Value = 2840
$ws.Range("L63").Value = 2900
$ws.Range("M63").Value = -2154
$ws.Range("N63").Value = -4272
$ws.Range("H66").Value = 2850
$ws.Range("I66").Value = 2840
$ws.Range("J66").Value = 2900
$ws.Range("K66").Value = 14200
$ws.Range("L66").Value = 14500
$ws.Range("M66").Value = -10768
$ws.Range("N66").Value = -21364
$ws.Range("I97").Value = 5682835.5
$ws.Range("J97").Value = 762.2
$ws.Range("K97").Value = 5682835.5
$ws.Range("L97").Value = 762.2
$ws.Range("M97").Value = -5682339.5
$ws.Range("N97").Value = -1754.2
$ws.Range("H132").Value = 45621.598
$ws.Range("I132").Value = 31109.303
$ws.Range("J132").Value = 79829.14
$ws.Range("K132").Value = 93327.909
$ws.Range("L132").Value = 239487.42
$ws.Range("M132").Value = -90797.909
$ws.Range("N132").Value = -244547.42
$ws.Range("H136").Value = 52737764
$ws.Range("I136").Value = 62563348
$ws.Range("K136").Value = 187690044
$ws.Range("M136").Value = -187687494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8196.308000000001
$ws.Range("I86").Value = 11016.6
$ws.Range("J86").Value = 3160.0715
$ws.Range("K86").Value = 11016.6
$ws.Range("L86").Value = 3160.0715
$ws.Range("M86").Value = -9893.6
$ws.Range("N86").Value = -5406.0715
$ws.Range("H89").Value = 8196.308000000001
$ws.Range("I89").Value = 11016.6
$ws.Range("J89").Value = 3160.0715
$ws.Range("K89").Value = 55083
$ws.Range("L89").Value = 15800.3575
$ws.Range("M89").Value = -49467
$ws.Range("N89").Value = -27032.3575
$ws.Range("H94").Value = 811.8889
$ws.Range("I94").Value = 819.5
$ws.Range("J94").Value = 796.6667
$ws.Range("K94").Value = 819.5
$ws.Range("L94").Value = 796.6667
$ws.Range("M94").Value = -368.5
$ws.Range("N94").Value = -1698.6667
$ws.Range("H99").Value = 1088.4849
$ws.Range("I99").Value = 1112.9412
$ws.Range("K99").Value = 1112.9412
$ws.Range("M99").Value = 385.0588
$ws.Range("H105").Value = 31252136
$ws.Range("I105").Value = 45456664
$ws.Range("J105").Value = 2180
$ws.Range("K105").Value = 45456664
$ws.Range("L105").Value = 2180
$ws.Range("M105").Value = -45454917
$ws.Range("N105").Value = -5674
$ws.Range("H134").Value = 4252.469
$ws.Range("I134").Value = 4128.4473
$ws.Range("J134").Value = 4680.909
$ws.Range("K134").Value = 12385.3419
$ws.Range("L134").Value = 14042.727
$ws.Range("M134").Value = -9850.341899999999
$ws.Range("N134").Value = -19112.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2988.244
$ws.Range("I31").Value = 1841.6666
$ws.Range("K31").Value = 1841.6666
$ws.Range("M31").Value = -1546.6666
$ws.Range("H34").Value = 2988.244
$ws.Range("I34").Value = 1841.6666
$ws.Range("K34").Value = 1841.6666
$ws.Range("M34").Value = -1639.6666
$ws.Range("H53").Value = 30600
$ws.Range("J53").Value = 30600
$ws.Range("L53").Value = 30600
$ws.Range("N53").Value = -31814
$ws.Range("H58").Value = 45456690
$ws.Range("I58").Value = 83335210
$ws.Range("J58").Value = 2461.3
$ws.Range("K58").Value = 83335210
$ws.Range("L58").Value = 2461.3
$ws.Range("M58").Value = -83335007
$ws.Range("N58").Value = -2867.3
$ws.Range("H94").Value = 8730
$ws.Range("I94").Value = 45000
$ws.Range("J94").Value = 1476
$ws.Range("K94").Value = 45000
$ws.Range("L94").Value = 1476
$ws.Range("M94").Value = -44549
$ws.Range("N94").Value = -2378
$ws.Range("H105").Value = 975.5625
$ws.Range("I105").Value = 932
$ws.Range("J105").Value = 1280.5
$ws.Range("K105").Value = 932
$ws.Range("L105").Value = 1280.5
$ws.Range("M105").Value = 815
$ws.Range("N105").Value = -4774.5
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H132").Value = 14635.173
$ws.Range("I132").Value = 1139.4445
$ws.Range("J132").Value = 49338.477
$ws.Range("K132").Value = 3418.3335
$ws.Range("L132").Value = 148015.431
$ws.Range("M132").Value = -888.3335000000002
$ws.Range("N132").Value = -153075.431
$ws.Range("H134").Value = 28266.627
$ws.Range("I134").Value = 2518
$ws.Range("K134").Value = 7554
$ws.Range("M134").Value = -5019
$ws.Range("H136").Value = 45456690
$ws.Range("I136").Value = 83335210
$ws.Range("J136").Value = 2461.3
$ws.Range("K136").Value = 250005630
$ws.Range("L136").Value = 7383.900000000001
$ws.Range("M136").Value = -250003080
$ws.Range("N136").Value = -12483.9
$ws.Range("H138").Value = 36570
$ws.Range("J138").Value = 36570
$ws.Range("L138").Value = 36570
$ws.Range("N138").Value = -46850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2186.4614
$ws.Range("I70").Value = 1060.5714
$ws.Range("J70").Value = 3500
$ws.Range("K70").Value = 3181.7142
$ws.Range("L70").Value = 10500
$ws.Range("M70").Value = -2866.7142
$ws.Range("N70").Value = -11130
$ws.Range("H73").Value = 2186.4614
$ws.Range("I73").Value = 1060.5714
$ws.Range("J73").Value = 3500
$ws.Range("K73").Value = 3181.7142
$ws.Range("L73").Value = 10500
$ws.Range("M73").Value = -2089.7142
$ws.Range("N73").Value = -12684
$ws.Range("H75").Value = 2666.6667
$ws.Range("H78").Value = 2666.6667
$ws.Range("H131").Value = 10417738
$ws.Range("J131").Value = 1117.1333
$ws.Range("L131").Value = 3351.3999
$ws.Range("N131").Value = -13431.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2321.2
$ws.Range("I126").Value = 2263.3845
$ws.Range("J126").Value = 2428.5715
$ws.Range("K126").Value = 6790.1535
$ws.Range("L126").Value = 7285.7145
$ws.Range("M126").Value = -4320.1535
$ws.Range("N126").Value = -12225.7145
$ws.Range("H132").Value = 97654.766
$ws.Range("I132").Value = 64671.125
$ws.Range("J132").Value = 203202.4
$ws.Range("K132").Value = 194013.375
$ws.Range("L132").Value = 609607.2
$ws.Range("M132").Value = -191483.375
$ws.Range("N132").Value = -614667.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2689.4736
$ws.Range("I16").Value = 1279.24
$ws.Range("J16").Value = 5401.4614
$ws.Range("K16").Value = 1279.24
$ws.Range("L16").Value = 5401.4614
$ws.Range("M16").Value = -1109.24
$ws.Range("N16").Value = -5741.4614
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -22246
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -71232
$ws.Range("H132").Value = 25402.766
$ws.Range("I132").Value = 11745.88
$ws.Range("J132").Value = 74177.36
$ws.Range("K132").Value = 35237.64
$ws.Range("L132").Value = 222532.08
$ws.Range("M132").Value = -32707.64
$ws.Range("N132").Value = -227592.08
$ws.Range("H136").Value = 52506.35
$ws.Range("I136").Value = 33544.195
$ws.Range("J136").Value = 117820.445
$ws.Range("K136").Value = 100632.585
$ws.Range("L136").Value = 353461.335
$ws.Range("M136").Value = -98082.58499999999
$ws.Range("N136").Value = -358561.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 92569.55
$ws.Range("I132").Value = 63890.438
$ws.Range("J132").Value = 169047.17
$ws.Range("K132").Value = 191671.314
$ws.Range("L132").Value = 507141.51
$ws.Range("M132").Value = -189141.314
$ws.Range("N132").Value = -512201.51
$ws.Range("H136").Value = 51659.023
$ws.Range("I136").Value = 50969.1
$ws.Range("K136").Value = 152907.3
$ws.Range("M136").Value = -150357.3
